# Updating End to End Suite
# Add a new worksheet "TC005" at the end of the workbook (after "TC004"),
# populate it with the menu/product header row + one data row, and make
# it the active/selected sheet (matching the xl/workbook.xml +
# xl/worksheets/sheet4.xml additions in the target diff).

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet so it lands at the end
# of the tab strip (TC002, TC003, TC004, TC005) and becomes the active tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "TC005"

# Header row
$ws.Range("A1").Value = "MenuName"
$ws.Range("B1").Value = "Product Id"
$ws.Range("C1").Value = "Quantity"
$ws.Range("D1").Value = "Size"
$ws.Range("E1").Value = "Color"

# Data row
$ws.Range("A2").Value = "Dresses"
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = "L"
$ws.Range("E2").Value = "Blue"

# Match the new sheet's saved selection (A1:E2) and keep it the active sheet.
$ws.Activate()
$ws.Range("A1:E2").Select()
